# Add two new lines to the "Cronologia Hernan Hamra" timeline document:
#   1. A "Running/entren." training-line, inserted right after the
#      "Personal trainer (23-40)" line (PARTE 2 block).
#   2. A "Raul (padre)" death-date line, inserted right after the
#      "Nace Matias (38)" line (PARTE 2 block).
#
# Both new paragraphs must match the surrounding paragraph formatting
# (w:spacing w:after="0" w:before="0"), which InsertParagraphAfter()
# naturally inherits from the paragraph it is invoked on.

$d = $word.ActiveDocument

$runningText = "└─ Running/entren. (23-40)│   │   │   ██▌ ██▌ ██▌ ██▌ ██▌ ██▌ ██▌ ██▌ ██▌ ██▌ ██▌ ██▌ ██▌ ██▌ ██▌ ██▌ ██▌ ██▌  2-3 media maratones"
$raulText    = "† Raúl (padre) (39)       │   │   │   │   │   │   │   │   │   │   │   │   │   │   │   │   │   │   │   ██▌ │  (12/06/2012) mismo día media maratón"

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*└─ Personal trainer (23-40)*sigue hasta 2022*") {
        $p.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs.Item($i + 1)
        $newPara.Range.Text = $runningText
        break
    }
}

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Nace Matías (38)*03/10*") {
        $p.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs.Item($i + 1)
        $newPara.Range.Text = $raulText
        break
    }
}

Write-Output "Done"
